$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$r = 96
$ws.Cells.Item($r, 1).Value = 7
$ws.Cells.Item($r, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item($r, 3).Value = "Ñuble"
$ws.Cells.Item($r, 4).Value = 45239
$ws.Cells.Item($r, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item($r, 5).Value = 16
$ws.Cells.Item($r, 6).Value = 100112026
$ws.Cells.Item($r, 7).Value = "Haba"
$ws.Cells.Item($r, 8).Value = "Sin especificar"
$ws.Cells.Item($r, 9).Value = "Primera"
$ws.Cells.Item($r, 10).Value = 80
$ws.Cells.Item($r, 11).Value = 10000
$ws.Cells.Item($r, 12).Value = 10000
$ws.Cells.Item($r, 13).Value = 10000
$ws.Cells.Item($r, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item($r, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item($r, 16).Value = 400
$ws.Cells.Item($r, 17).Value = 25
$ws.Cells.Item($r, 18).Value = "Hortaliza"
